$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IDXX")

# Row 4 - Inventory
$ws.Range("B4").Value = 210000000.0
$ws.Range("C4").Value = 218000000.0
$ws.Range("D4").Value = 229000000.0
$ws.Range("E4").Value = 212000000.0
$ws.Range("F4").Value = 195000000.0

# Row 14 - Accounts Payable
$ws.Range("B14").Value = 75000000.0
$ws.Range("C14").Value = 72000000.0
$ws.Range("D14").Value = 74000000.0
$ws.Range("E14").Value = 74000000.0
$ws.Range("F14").Value = 72000000.0

# Row 22 - Long Term Tax Liability (Deferred)
$ws.Range("B22").Value = -20000000.0
$ws.Range("C22").Value = 19000000.0
$ws.Range("D22").Value = 29000000.0
$ws.Range("E22").Value = 32000000.0
$ws.Range("F22").Value = 25000000.0
